$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 220609
$ws.Range("B4").Value = "박영서"
$ws.Range("C4").Value = "VGG16"
$ws.Range("D4").Value = 11
$ws.Range("E4").Value = 0.8494
$ws.Range("F4").Value = 0.8415
$ws.Range("G4").Value = 0.655
$ws.Range("H4").Value = 0.7116

$ws.Range("H7").Select()
